# Fruta / hortaliza, semanal
# Updates the weekly price data (date, volume, min/max/avg price, price per kg)
# for rows 2-19, keeping other constant columns (A,B,C,E,F,G,H,I,N,O,Q,R) unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)
$rows = @{
    2  = @(44677, 20,  5500, 5500, 5500, 5500)
    3  = @(44685, 60,  5000, 6000, 5333, 5333)
    4  = @(44476, 30,  2200, 2200, 2200, 2200)
    5  = @(44707, 100, 4700, 4700, 4700, 4700)
    7  = @(44706, 90,  4700, 4700, 4700, 4700)
    8  = @(44473, 140, 1600, 1600, 1600, 1600)
    9  = @(44669, 60,  6250, 6250, 6250, 6250)
    11 = @(44452, 120, 2300, 2300, 2300, 2300)
    12 = @(44447, 75,  2200, 2200, 2200, 2200)
    13 = @(44496, 40,  2200, 2200, 2200, 2200)
    14 = @(44203, 30,  2000, 2000, 2000, 2000)
    15 = @(44679, 30,  5500, 5500, 5500, 5500)
    16 = @(44497, 50,  2200, 2200, 2200, 2200)
    17 = @(44453, 20,  2300, 2300, 2300, 2300)
    18 = @(44474, 20,  1600, 1600, 1600, 1600)
    19 = @(44487, 50,  2200, 2200, 2200, 2200)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]

    $ws.Cells.Item($r, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $vals[1]   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $vals[2]   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals[3]   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals[4]   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $vals[5]   # P - Precio $/Kg
}
